$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they remain text like the source data.
$textFormatCells = @("D5","D6","D7","D8","D11","D14","D16","D19","D20","D21","D22","D23","D24","D25","D28","D31","D33","D34","D36","D37","D38","D39","D41","D42","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values row by row as described by the diff
# Row 2
$ws.Cells.Item(2, 4).Value = '88.912.10'
$ws.Cells.Item(2, 5).Value = '  +3.59%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '3.283.00'
$ws.Cells.Item(3, 5).Value = '  -1.40%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
# Row 5
$ws.Cells.Item(5, 4).Value = '213.37'
$ws.Cells.Item(5, 5).Value = '  -2.70%  '
# Row 6
$ws.Cells.Item(6, 4).Value = '630.56'
$ws.Cells.Item(6, 5).Value = '  -1.25%  '
# Row 7
$ws.Cells.Item(7, 4).Value = '0.391'
$ws.Cells.Item(7, 5).Value = '  +20.83%  '
# Row 8
$ws.Cells.Item(8, 4).Value = '0.693'
$ws.Cells.Item(8, 5).Value = '  +16.51%  '
# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.07%  '
# Row 10
$ws.Cells.Item(10, 4).Value = '3.280.27'
$ws.Cells.Item(10, 5).Value = '  -1.28%  '
# Row 11
$ws.Cells.Item(11, 4).Value = '0.579'
$ws.Cells.Item(11, 5).Value = '  -3.12%  '
# Row 12
$ws.Cells.Item(12, 5).Value = '  +11.75%  '
# Row 13
$ws.Cells.Item(13, 5).Value = '  -4.40%  '
# Row 14
$ws.Cells.Item(14, 4).Value = '34.13'
$ws.Cells.Item(14, 5).Value = '  -0.45%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '3.887.33'
$ws.Cells.Item(15, 5).Value = '  -1.05%  '
# Row 16
$ws.Cells.Item(16, 4).Value = '5.38'
$ws.Cells.Item(16, 5).Value = '  -0.64%  '
# Row 17
$ws.Cells.Item(17, 4).Value = '88.443.12'
$ws.Cells.Item(17, 5).Value = '  +3.58%  '
# Row 18
$ws.Cells.Item(18, 4).Value = '3.289.66'
$ws.Cells.Item(18, 5).Value = '  -0.46%  '
# Row 19
$ws.Cells.Item(19, 4).Value = '3.11'
$ws.Cells.Item(19, 5).Value = '  -2.16%  '
# Row 20
$ws.Cells.Item(20, 4).Value = '14.12'
$ws.Cells.Item(20, 5).Value = '  -3.91%  '
# Row 21
$ws.Cells.Item(21, 4).Value = '436.94'
$ws.Cells.Item(21, 5).Value = '  -0.81%  '
# Row 22
$ws.Cells.Item(22, 4).Value = '8.90'
$ws.Cells.Item(22, 5).Value = '  -3.11%  '
# Row 23
$ws.Cells.Item(23, 4).Value = '5.39'
$ws.Cells.Item(23, 5).Value = '  +2.48%  '
# Row 24
$ws.Cells.Item(24, 4).Value = '7.39'
$ws.Cells.Item(24, 5).Value = '  +0.14%  '
# Row 25
$ws.Cells.Item(25, 4).Value = '12.32'
$ws.Cells.Item(25, 5).Value = '  +0.67%  '
# Row 26
$ws.Cells.Item(26, 5).Value = '  -5.44%  '
# Row 27
$ws.Cells.Item(27, 4).Value = '3.442.35'
$ws.Cells.Item(27, 5).Value = '  -1.15%  '
# Row 28
$ws.Cells.Item(28, 4).Value = '77.06'
$ws.Cells.Item(28, 5).Value = '  -1.70%  '
# Row 29
$ws.Cells.Item(29, 5).Value = '  +2.65%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.01%  '
# Row 31
$ws.Cells.Item(31, 4).Value = '0.190'
$ws.Cells.Item(31, 5).Value = '  +13.20%  '
# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.40%  '
# Row 33
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 4).Value = '8.90'
$ws.Cells.Item(33, 5).Value = '  -4.05%  '
# Row 34
$ws.Cells.Item(34, 2).Value = 'Bittensor'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(34, 4).Value = '568.22'
$ws.Cells.Item(34, 5).Value = '  -8.40%  '
# Row 35
$ws.Cells.Item(35, 5).Value = '  -10.11%  '
# Row 36
$ws.Cells.Item(36, 2).Value = 'RenderToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(36, 4).Value = '7.20'
$ws.Cells.Item(36, 5).Value = '  +11.07%  '
# Row 37
$ws.Cells.Item(37, 2).Value = 'PancakeSwap'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(37, 4).Value = '1.97'
$ws.Cells.Item(37, 5).Value = '  -4.10%  '
# Row 38
$ws.Cells.Item(38, 4).Value = '0.138'
$ws.Cells.Item(38, 5).Value = '  -8.87%  '
# Row 39
$ws.Cells.Item(39, 4).Value = '22.66'
$ws.Cells.Item(39, 5).Value = '  -2.77%  '
# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.40%  '
# Row 41
$ws.Cells.Item(41, 4).Value = '21.80'
$ws.Cells.Item(41, 5).Value = '  +2.49%  '
# Row 42
$ws.Cells.Item(42, 4).Value = '0.399'
$ws.Cells.Item(42, 5).Value = '  -4.97%  '
# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.97%  '
# Row 44
$ws.Cells.Item(44, 4).Value = '3.00'
$ws.Cells.Item(44, 5).Value = '  -3.04%  '
# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.11%  '
# Row 46
$ws.Cells.Item(46, 4).Value = '154.19'
$ws.Cells.Item(46, 5).Value = '  -3.08%  '
# Row 47
$ws.Cells.Item(47, 4).Value = '180.86'
$ws.Cells.Item(47, 5).Value = '  -4.88%  '
# Row 48
$ws.Cells.Item(48, 4).Value = '45.07'
$ws.Cells.Item(48, 5).Value = '  -0.15%  '
# Row 49
$ws.Cells.Item(49, 4).Value = '1.30'
$ws.Cells.Item(49, 5).Value = '  -5.31%  '
# Row 50
$ws.Cells.Item(50, 2).Value = 'Filecoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(50, 4).Value = '4.25'
$ws.Cells.Item(50, 5).Value = '  -0.69%  '
# Row 51
$ws.Cells.Item(51, 4).Value = '0.0678'
$ws.Cells.Item(51, 5).Value = '  +20.38%  '

"done"